# Full resnet-50 graph is completed for both inference and BP.
# Insert a new opcode row "MULT_DER" right after the existing "MULT" row
# (row 7), pushing ADD / ACT / POOL / POOL_BP / BNORM / BNORM_BP / NOP / EOL
# down by one row, and decrement the "# of operands" column for every one
# of those opcodes by 1 (MULT_DER absorbs one of the operands that used to
# be counted against the following opcode).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 7 (MULT) into a freshly inserted row 8 so the new row
# naturally inherits the same cell styles/formatting as the rest of the
# opcode table, then overwrite its contents for MULT_DER.
$ws.Rows("7:7").Copy()
$ws.Rows("8:8").Insert()

# New row 8: MULT_DER, 1 operand, "no of iteration"
$ws.Range("A8").Value = "MULT_DER"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "no of iteration"

# MULT (row 7) now has only 1 operand as well.
$ws.Range("B7").Value = 1

# Rows 9-15 are the old ADD(8)..NOP(14) rows shifted down by one; each of
# their operand counts drops by 1 (EOL/row16 keeps its 0).
$ws.Range("B9").Value = 1    # ADD (was 2)
$ws.Range("B10").Value = 2   # ACT (was 3)
$ws.Range("B11").Value = 2   # POOL (was 3)
$ws.Range("B12").Value = 2   # POOL_BP (was 3)
$ws.Range("B13").Value = 2   # BNORM (was 3)
$ws.Range("B14").Value = 2   # BNORM_BP (was 3)
$ws.Range("B15").Value = 1   # NOP (was 2)

# Leave the cursor where the author ended up editing: the new EOL row.
$ws.Range("B16").Select() | Out-Null
